$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of invoice data (row 3) mirroring the existing row 2 layout:
#   A: id (number)       B: title (text)   C: count (number)
#   D: price (text, numeric-looking but stored as a string)
#   E: created_at (date/time serial, same display format as E2)
#   F: total_price (text, numeric-looking but stored as a string)

$ws.Range("A3").Value = 45
$ws.Range("B3").Value = "Product 2"
$ws.Range("C3").Value = 6

# D3 and F3 must hold the literal text "32000.00" / "192000.00" (not be
# auto-converted to numbers) while keeping the default/General style - the
# same way the existing D2/F2 text cells look. Typing straight into a
# General-formatted cell makes Excel coerce numeric-looking text to a
# number, and pre-formatting the cell as Text mints a brand-new style.
# Instead, stage the literal string as a text formula result on a scratch
# cell and paste-special just the values in - Excel keeps the source's
# text type without touching number formats/styles.
$ws.Range("H1").Formula = "=""32000.00"""
$ws.Range("H1").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("H1").Formula = "=""192000.00"""
$ws.Range("H1").Copy()
$ws.Range("F3").PasteSpecial(-4163)

$ws.Range("H1").ClearContents()

# E3: same created_at style as E2 (reuse the existing YYYY-MM-DD HH:MM:SS
# number format instead of minting a new one).
$ws.Range("E3").Value = 45687.78773686659
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
